$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.011618333333334
$ws.Range("H2").Value = 3.034855
$ws.Range("I2").Value = 0.5235149663433657
$ws.Range("J2").Value = 0.5235149663433657
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 14.48827738741834
$ws.Range("R2").Value = 130.394496486765
$ws.Range("S2").Value = 0.1544143608311015
$ws.Range("T2").Value = 0.1544143608311015

$ws.Range("G3").Value = 1.011618333333334
$ws.Range("H3").Value = 3.034855
$ws.Range("I3").Value = 0.5235149663433657
$ws.Range("J3").Value = 0.5235149663433657
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 27.39884190349834
$ws.Range("R3").Value = 246.589577131485
$ws.Range("S3").Value = 0.2920136429548978
$ws.Range("T3").Value = 0.2920136429548978

$ws.Range("G4").Value = 1.011618333333334
$ws.Range("H4").Value = 3.034855
$ws.Range("I4").Value = 0.5235149663433657
$ws.Range("J4").Value = 0.5235149663433657
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 7.232858980689445
$ws.Range("R4").Value = 65.095730826205
$ws.Range("S4").Value = 0.07708696255736626
$ws.Range("T4").Value = 0.07708696255736625

$ws.Range("I5").Value = 0.2899264353016711
$ws.Range("J5").Value = 0.2899264353016712
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 8.023714481241665
$ws.Range("R5").Value = 72.213430331175
$ws.Range("S5").Value = 0.08551580771003989
$ws.Range("T5").Value = 0.0855158077100399

$ws.Range("I6").Value = 0.2899264353016711
$ws.Range("J6").Value = 0.2899264353016712
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("S6").Value = 0.1617193012698697
$ws.Range("T6").Value = 0.1617193012698698

$ws.Range("I7").Value = 0.2899264353016711
$ws.Range("J7").Value = 0.2899264353016712
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 4.005610452663888
$ws.Range("R7").Value = 36.050494073975
$ws.Range("S7").Value = 0.04269132632176146
$ws.Range("T7").Value = 0.04269132632176146

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.360498
$ws.Range("H8").Value = 1.081494
$ws.Range("I8").Value = 0.1865585983549632
$ws.Range("J8").Value = 0.1865585983549632
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 5.163009456737999
$ws.Range("R8").Value = 46.467085110642
$ws.Range("S8").Value = 0.05502674913716513
$ws.Range("T8").Value = 0.05502674913716513

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.360498
$ws.Range("H9").Value = 1.081494
$ws.Range("I9").Value = 0.1865585983549632
$ws.Range("J9").Value = 0.1865585983549632
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 9.763788756162
$ws.Range("R9").Value = 87.874098805458
$ws.Range("S9").Value = 0.1040613152107314
$ws.Range("T9").Value = 0.1040613152107314

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.360498
$ws.Range("H10").Value = 1.081494
$ws.Range("I10").Value = 0.1865585983549632
$ws.Range("J10").Value = 0.1865585983549632
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 2.577485115586
$ws.Range("R10").Value = 23.197366040274
$ws.Range("S10").Value = 0.02747053400706665
$ws.Range("T10").Value = 0.02747053400706664
